$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Revision de no conformidades: mark rows 5 and 6 (STATUS column F) as "Cerrada"
# instead of "En proceso". This makes the shared string "En proceso" (previously
# duplicated) become unused in the string table once these are the only rows
# referencing it, and Excel will drop the now-unused entry on save.
$ws.Range("F5").Value = "Cerrada"
$ws.Range("F6").Value = "Cerrada"

# Update the active cell/selection to F6 to match the saved view state.
$ws.Range("F6").Select()
